# Append a new data row (row 45) with the latest reading pulled from Adafruit IO,
# matching the existing table's layout: Timestamp, Feed Key, Value, Latitude, Longitude, Elevation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 45

$ws.Range("A$newRow").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$newRow").Value = "temperature"

# "25" looks numeric, but the rest of the column stores it as text (like the other
# rows), so use the leading-apostrophe trick to force text, then re-apply the plain
# (unstyled) format from the row above so no stray cell style is introduced.
$ws.Range("C$newRow").Value = "'25"
$ws.Range("C$newRow").Style = $ws.Range("C44").Style

$ws.Range("D$newRow").Value = "N/A"
$ws.Range("E$newRow").Value = "N/A"
$ws.Range("F$newRow").Value = "N/A"
